# Apply updated crypto price / volume data to the worksheet.
# Numeric-looking Price values are prefixed with a literal apostrophe so Excel
# stores them as text (matching the original inlineStr cells) instead of
# silently converting them to floating point numbers (which would drop
# trailing zeros / change formatting).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.359.46'
$ws.Range("E2").Value = '  +0.81%  '
$ws.Range("D3").Value = '1.688.20'
$ws.Range("E3").Value = '  +0.89%  '
$ws.Range("E4").Value = '  +0.75%  '
$ws.Range("D5").Value = '''218.43'
$ws.Range("E5").Value = '  +0.67%  '
$ws.Range("D6").Value = '''0.5466'
$ws.Range("E6").Value = '  +4.70%  '
$ws.Range("E7").Value = '  +0.73%  '
$ws.Range("D8").Value = '''0.2728'
$ws.Range("E8").Value = '  +1.01%  '
$ws.Range("E9").Value = '  +1.02%  '
$ws.Range("E10").Value = '  +0.87%  '
$ws.Range("D11").Value = '''0.07679'
$ws.Range("E11").Value = '  +3.29%  '
$ws.Range("D12").Value = '1.699.65'
$ws.Range("E12").Value = '  +1.38%  '
$ws.Range("D13").Value = '''4.537'
$ws.Range("E13").Value = '  +0.34%  '
$ws.Range("D14").Value = '''0.5811'
$ws.Range("E14").Value = '  -0.39%  '
$ws.Range("D15").Value = '''0.000008330'
$ws.Range("E15").Value = '  -2.37%  '
$ws.Range("D16").Value = '''65.13'
$ws.Range("E16").Value = '  +1.36%  '
$ws.Range("D17").Value = '26.411.53'
$ws.Range("E17").Value = '  +1.80%  '
$ws.Range("E18").Value = '  +0.21%  '
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("E20").Value = '  +1.57%  '
$ws.Range("D21").Value = '''190.66'
$ws.Range("E21").Value = '  +0.37%  '
$ws.Range("D22").Value = '''6.230'
$ws.Range("E22").Value = '  +0.63%  '
$ws.Range("D23").Value = '''1.011'
$ws.Range("E23").Value = '  +0.61%  '
$ws.Range("D24").Value = '''149.40'
$ws.Range("E24").Value = '  +3.12%  '
$ws.Range("E25").Value = '  +5.03%  '
$ws.Range("D26").Value = '''7.888'
$ws.Range("E26").Value = '  +3.53%  '
$ws.Range("D27").Value = '''15.70'
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("D28").Value = '''0.06359'
$ws.Range("E28").Value = '  -3.28%  '
$ws.Range("D29").Value = '''1.410'
$ws.Range("E29").Value = '  +5.50%  '
$ws.Range("E30").Value = '  +0.74%  '
$ws.Range("D31").Value = '''3.576'
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("D32").Value = '''3.578'
$ws.Range("E32").Value = '  +1.28%  '
$ws.Range("E33").Value = '  +0.36%  '
$ws.Range("D34").Value = '''1.043'
$ws.Range("E34").Value = '  +2.28%  '
$ws.Range("D35").Value = '''0.6210'
$ws.Range("E35").Value = '  +0.69%  '
$ws.Range("D36").Value = '''2.414'
$ws.Range("E36").Value = '  +2.03%  '
$ws.Range("D37").Value = '''2.720'
$ws.Range("E37").Value = '  +0.74%  '
$ws.Range("D38").Value = '''6.238'
$ws.Range("E38").Value = '  -0.66%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '1.117.96'
$ws.Range("E39").Value = '  +2.11%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '''0.01630'
$ws.Range("E40").Value = '  +1.98%  '
$ws.Range("D41").Value = '''0.8793'
$ws.Range("E41").Value = '  +0.87%  '
$ws.Range("E42").Value = '  +0.40%  '
$ws.Range("D43").Value = '''100.96'
$ws.Range("E43").Value = '  +0.12%  '
$ws.Range("D44").Value = '1.840.56'
$ws.Range("E44").Value = '  +1.17%  '
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").Value = '''57.31'
$ws.Range("E46").Value = '  +1.40%  '
$ws.Range("D47").Value = '''8.215'
$ws.Range("E47").Value = '  +0.68%  '
$ws.Range("E48").Value = '  +1.23%  '
$ws.Range("D49").Value = '''0.05271'
$ws.Range("E49").Value = '  +0.60%  '
$ws.Range("D50").Value = '''0.4305'
$ws.Range("E50").Value = '  +0.59%  '
$ws.Range("D51").Value = '''6.037'
$ws.Range("E51").Value = '  +0.70%  '
